$d = $word.ActiveDocument
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Insert two new paragraphs right before the "Proces" Heading1 paragraph:
#      - an empty (Normal) paragraph
#      - a Normal paragraph with "Startet med DDS lite, gået over til EF",
#        carrying the _GoBack bookmark that used to sit on the
#        "Cultureinfo..." paragraph.
#    A trailing dummy paragraph is appended to the inserted XML so that the
#    final (merging) paragraph mark doesn't fuse into "Proces" and steal its
#    Heading1 style; the dummy is deleted again right after.
# ---------------------------------------------------------------------------

$r = $d.Content
$null = $r.Find.Execute("Proces")
$procesIdx = $r.Paragraphs.Item(1).Index
$prevPara = $d.Paragraphs.Item($procesIdx - 1)
$insPos = $prevPara.Range.End - 1

$newBlockXml = "<w:p $xmlNs><w:pPr><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr></w:pPr></w:p>" + `
  "<w:p $xmlNs><w:pPr><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr><w:t>Startet med DDS lite, gået over til EF</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>" + `
  "<w:p $xmlNs><w:pPr><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr></w:pPr></w:p>"

$insRange = $d.Range($insPos, $insPos)
$insRange.InsertXML($newBlockXml)

$r2 = $d.Content
$null = $r2.Find.Execute("Proces")
$procesIdx2 = $r2.Paragraphs.Item(1).Index
$dummyPara = $d.Paragraphs.Item($procesIdx2 - 1)
$dummyPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Strip the _GoBack bookmark from the "Cultureinfo..." paragraph (it moved
#    to the new paragraph above) by deleting and rewriting the paragraph.
#    Same trailing-dummy trick as above to keep paragraph boundaries intact.
# ---------------------------------------------------------------------------

$cr = $d.Content
$null = $cr.Find.Execute("Cultureinfo skal fixes")
$cultIdx = $cr.Paragraphs.Item(1).Index
$cultPara = $d.Paragraphs.Item($cultIdx)
$cultStart = $cultPara.Range.Start
$cultFull = $d.Range($cultStart, $cultPara.Range.End)
$cultFull.Delete()

$cultXml = "<w:p $xmlNs><w:pPr><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr><w:t>Cultureinfo skal fixes – Flere muligheder. Set på host PC eller tving gennem kode.</w:t></w:r></w:p>" + `
  "<w:p $xmlNs><w:pPr><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr></w:pPr></w:p>"
$cultInsRange = $d.Range($cultStart, $cultStart)
$cultInsRange.InsertXML($cultXml)

$sr = $d.Content
$null = $sr.Find.Execute("Som systemet er designet nu")
$somIdx = $sr.Paragraphs.Item(1).Index
$dummyPara2 = $d.Paragraphs.Item($somIdx - 1)
$dummyPara2.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Remove the lastRenderedPageBreak marker from the "Som systemet..."
#    paragraph by rewriting its content in place.
# ---------------------------------------------------------------------------

$sr2 = $d.Content
$null = $sr2.Find.Execute("Som systemet er designet nu")
$somIdx2 = $sr2.Paragraphs.Item(1).Index
$somPara = $d.Paragraphs.Item($somIdx2)
$somFullRange = $d.Range($somPara.Range.Start, $somPara.Range.End)
$somXml = "<w:p $xmlNs><w:pPr><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"da-DK`"/></w:rPr><w:t>Som systemet er designet nu, lægger respektive sensordata for alle pools i deres enkelte tabel.</w:t></w:r></w:p>"
$somFullRange.InsertXML($somXml)

Write-Output "done"
